$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
  2  = @(2,2,1,2,2,2)
  3  = @(2,2,1,1,2,2)
  4  = @(2,2,1,2,2,2)
  5  = @(2,2,1,1,1,2)
  6  = @(2,2,2,1,1,2)
  7  = @(2,1,1,0,2,2)
  8  = @(2,2,1,1,2,2)
  9  = @(2,2,1,1,2,2)
  10 = @(2,1,1,2,1,2)
  11 = @(2,1,1,1,1,2)
  12 = @(2,0,1,1,1,2)
  13 = @(2,2,1,1,1,2)
  14 = @(2,1,1,1,1,2)
  15 = @(2,2,2,1,1,2)
  16 = @(2,2,2,1,2,2)
  17 = @(2,2,2,1,1,2)
  18 = @(2,1,1,1,1,2)
  19 = @(2,1,1,1,1,2)
  20 = @(2,1,1,1,1,2)
  21 = @(2,2,1,1,2,2)
  22 = @(2,2,1,1,1,2)
  23 = @(2,2,2,1,2,2)
  24 = @(2,2,1,1,1,2)
  25 = @(2,2,1,1,1,2)
  26 = @(2,2,2,2,1,2)
  27 = @(2,1,2,1,2,2)
  28 = @(2,2,1,2,1,2)
  29 = @(2,1,1,1,1,2)
  30 = @(2,2,1,1,1,2)
  31 = @(2,1,1,1,1,2)
  32 = @(2,1,1,1,1,2)
  33 = @(2,1,1,1,1,2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i  # column E = 5
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# Sheet view changes: zoom, frozen pane, selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 29
$ws.Range("E34").Select()
